# Cut 20h to 10h
# - Drop the stale "Лист2" summary sheet
# - Re-run several rows with a 36000s (10h) time budget instead of 72000s (20h),
#   updating the LastLB/Layers/LB columns to match the new run, and re-deriving
#   the Gap column (I) via formula where appropriate
# - Move the active selection to E21

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Лист1")

# --- Remove the old "Лист2" sheet -----------------------------------------
$ws2 = $wb.Worksheets.Item("Лист2")
$ws2.Delete()

# --- Row 6 (ESC47) ----------------------------------------------------------
$ws.Range("E6").Value = 36000
$ws.Range("F6").Value = 18020
$ws.Range("G6").Value = 6
$ws.Range("J6").Value = 980
$ws.Range("I6").Formula = "=(D6-J6)/J6"

# --- Row 16 (ft70.4) ---------------------------------------------------------
$ws.Range("I16").Formula = "=(D16-J16)/J16"
$ws.Range("Q16").Formula = "=(D16-R16)/R16"

# --- Row 17 (kro124p.1) ------------------------------------------------------
$ws.Range("E17").Value = 36000
$ws.Range("F17").Formula = "=1+18+348"
$ws.Range("G17").Value = 3
$ws.Range("J17").Value = 27869
$ws.Range("I17").Formula = "=(D17-J17)/J17"

# --- Row 18 (kro124p.2) ------------------------------------------------------
$ws.Range("E18").Value = 36000
$ws.Range("F18").Value = 25773
$ws.Range("I18").Formula = "=(D18-J18)/J18"

# --- Row 19 (kro124p.3) ------------------------------------------------------
$ws.Range("E19").Value = 36000
$ws.Range("F19").Value = 31863

# --- Row 20 (kro124p.4) ------------------------------------------------------
$ws.Range("E20").Value = 36000
$ws.Range("F20").Value = 11405
$ws.Range("G20").Value = 6
$ws.Range("J20").Value = 38137
$ws.Range("I20").Formula = "=(D20-J20)/J20"

# --- Row 34 (rbg341a) ---------------------------------------------------------
$ws.Range("E34").Value = 36000
$ws.Range("F34").Value = 16392
$ws.Range("G34").Value = 6

# --- Row 35 (rbg358a) ---------------------------------------------------------
$ws.Range("E35").Value = 36000
$ws.Range("F35").Value = 16757
$ws.Range("G35").Value = 8

# --- Row 36 (rbg378a) ---------------------------------------------------------
$ws.Range("E36").Value = 36000

# --- Restore the view: scroll back to the top and select E21 ---------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("E21").Select()
